$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
  $r = $ws.Range($addr)
  $r.NumberFormat = "@"
  $r.Value = $val
  $r.Style = "Normal"
}

Set-TextValue "D2" "30.051.46"
Set-TextValue "E2" "  +2.42%  "
Set-TextValue "D3" "1.887.75"
Set-TextValue "E3" "  +2.60%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.50%  "
Set-TextValue "D5" "245.67"
Set-TextValue "E5" "  -0.45%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.49%  "
Set-TextValue "D7" "0.4957"
Set-TextValue "E7" "  +0.36%  "
Set-TextValue "D8" "44.33"
Set-TextValue "E8" "  -0.31%  "
Set-TextValue "D9" "0.2915"
Set-TextValue "E9" "  +4.36%  "
Set-TextValue "D10" "0.06611"
Set-TextValue "E10" "  +3.15%  "
Set-TextValue "D11" "1.892.01"
Set-TextValue "E11" "  +2.90%  "
Set-TextValue "D12" "16.90"
Set-TextValue "E12" "  +0.69%  "
Set-TextValue "D13" "0.07199"
Set-TextValue "E13" "  +1.17%  "
Set-TextValue "D14" "0.6707"
Set-TextValue "E14" "  +2.84%  "
Set-TextValue "D15" "85.45"
Set-TextValue "E15" "  +1.28%  "
Set-TextValue "D16" "4.819"
Set-TextValue "E16" "  +2.17%  "
Set-TextValue "D17" "30.088.45"
Set-TextValue "E17" "  +2.60%  "
Set-TextValue "D18" "0.000007880"
Set-TextValue "E18" "  +7.25%  "
Set-TextValue "D19" "1.000"
Set-TextValue "E19" "  +0.24%  "
Set-TextValue "D20" "12.77"
Set-TextValue "E20" "  +3.66%  "
Set-TextValue "D21" "2.135.65"
Set-TextValue "E21" "  +3.69%  "
Set-TextValue "D22" "1.002"
Set-TextValue "E22" "  +0.51%  "
Set-TextValue "D23" "4.763"
Set-TextValue "E23" "  +4.30%  "
Set-TextValue "D24" "5.586"
Set-TextValue "E24" "  +2.88%  "
Set-TextValue "D25" "9.126"
Set-TextValue "E25" "  +2.59%  "
Set-TextValue "D26" "148.49"
Set-TextValue "E26" "  +3.25%  "
Set-TextValue "D27" "133.47"
Set-TextValue "E27" "  +0.42%  "
Set-TextValue "D28" "16.73"
Set-TextValue "E28" "  +1.54%  "
Set-TextValue "D29" "1.938"
Set-TextValue "E29" "  +1.85%  "
Set-TextValue "D30" "1.377"
Set-TextValue "E30" "  -1.62%  "
Set-TextValue "D31" "4.176"
Set-TextValue "E31" "  +0.77%  "
Set-TextValue "D32" "0.08698"
Set-TextValue "E32" "  +3.88%  "
Set-TextValue "D33" "3.933"
Set-TextValue "E33" "  +3.61%  "
Set-TextValue "D34" "0.05106"
Set-TextValue "E34" "  +3.36%  "
Set-TextValue "D35" "1.115"
Set-TextValue "E35" "  +0.92%  "
Set-TextValue "D36" "0.7043"
Set-TextValue "E36" "  +4.54%  "
Set-TextValue "D37" "2.673"
Set-TextValue "E37" "  -0.37%  "
Set-TextValue "D38" "2.214"
Set-TextValue "E38" "  -3.50%  "
Set-TextValue "D39" "2.703"
Set-TextValue "E39" "  -0.54%  "
Set-TextValue "D40" "0.9391"
Set-TextValue "E40" "  -1.38%  "
Set-TextValue "D41" "0.01647"
Set-TextValue "E41" "  +3.19%  "
Set-TextValue "D42" "6.063"
Set-TextValue "E42" "  -2.08%  "
Set-TextValue "D43" "0.9988"
Set-TextValue "E43" "  +0.09%  "
Set-TextValue "D44" "103.28"
Set-TextValue "E44" "  +0.84%  "
Set-TextValue "D45" "0.4179"
Set-TextValue "E45" "  +2.18%  "
Set-TextValue "D46" "7.462"
Set-TextValue "E46" "  +3.14%  "
Set-TextValue "D47" "0.1260"
Set-TextValue "E47" "  +3.07%  "
Set-TextValue "D48" "0.05725"
Set-TextValue "E48" "  +2.80%  "
Set-TextValue "D49" "32.66"
Set-TextValue "E49" "  +2.70%  "
Set-TextValue "D50" "8.206"
Set-TextValue "E50" "  +1.32%  "
Set-TextValue "D51" "0.3715"
Set-TextValue "E51" "  +2.62%  "
